# Daily auto-push edit: insert a new sampling row for 2026/02/20 (Friday)
# ahead of the existing "2026/12/29" block, shifting every subsequent row
# down by one (row 836 -> 837, ... 877 -> 878).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 836; everything currently at/after
# row 836 (through the old last row, 877) shifts down to 837..878.
$xlShiftDown = -4121
$ws.Rows.Item(836).Insert($xlShiftDown)

# Column A holds a "YYYY/MM/DD" text label (not a real date), so force the
# cell to Text format before writing it — otherwise Excel's autodetect
# would silently convert the literal into a date serial number.
$ws.Range("A836").NumberFormat = "@"
$ws.Range("A836").Value = "2026/02/20"
$ws.Range("B836").Value = "金"
$ws.Range("C836").Value = 13
$ws.Range("D836").Value = 201

# Drop the temporary text-format override so the new row's style matches
# its unstyled neighbors (no explicit cell style, same as the rest of the
# data rows in the sheet).
$ws.Range("A836").ClearFormats()
